# Edit script: updates slide numbering / titles for the "Define Query
# Variable: Part 3" section (renumbering 10.4/10.5/10.6 -> 10.3/10.4/10.5),
# and resizes/repositions a couple of shapes on slide 10 to match the
# taller subtitle text box.

$p = $ppt.ActivePresentation
$EMU_PER_PT = 12700

# ---------------------------------------------------------------------
# Slide 1: Title slide - "10 Query Variable: Part 3" -> "10 Define Query
# Variable: Part 3"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "10 Define Query Variable: Part 3"

# ---------------------------------------------------------------------
# Slide 2: Agenda/overview slide - same title text change
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "10 Define Query Variable: Part 3"

# ---------------------------------------------------------------------
# Slide 8: "10.4 Run Application" -> "10.3 Step 3" (title) and
# "Step 4: Run Application" -> "Set 3: Define Variable" (subtitle first run)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "10.3 Step 3"
$s8sub = $s8.Shapes.Item(2).TextFrame.TextRange
$s8p1 = $s8sub.Paragraphs(1,1)
$s8sub.Characters($s8p1.Start, $s8p1.Length).Text = "Set 3: Define Variable"

# ---------------------------------------------------------------------
# Slide 9: "10.5 Run Application" -> "10.4 Run Application" (title)
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "10.4 Run Application"

# ---------------------------------------------------------------------
# Slide 10: "10.5 Run Application" -> "10.4 Run Application" (title),
# subtitle placeholder grows taller, and the screenshot picture below it
# moves down to match.
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "10.4 Run Application"

$s10sub = $s10.Shapes.Item(2)
$s10sub.Height = (1343176 / $EMU_PER_PT) + 0.00002

$s10pic = $s10.Shapes.Item(6)
$s10pic.Top = (2915726 / $EMU_PER_PT) + 0.00002

# ---------------------------------------------------------------------
# Slide 12: "10.6 Verify" -> "10.5 Verify" (title)
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(2).TextFrame.TextRange.Text = "10.5 Verify"

# ---------------------------------------------------------------------
# Slide 13: "10.6 Verify" -> "10.5 Verify" (title)
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "10.5 Verify"
